$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the cell updates described by the diff (updated coin prices/ranks
# and the "Worstin24h"/"Bestin24h" label reshuffle from the GitHub Actions run).
# Numeric-looking values are written with a leading apostrophe so Excel stores
# them as literal text (matching the original inlineStr/text cells) instead of
# converting them to real numbers, which would drop formatting such as trailing zeros.

$ws.Range("D2").Value = "'246.56"
$ws.Range("G2").Value = "'18"
$ws.Range("D3").Value = "'26.37"
$ws.Range("G3").Value = "'18"
$ws.Range("D4").Value = "'5.086"
$ws.Range("G4").Value = "'18"
$ws.Range("G5").Value = "'18"
$ws.Range("D6").Value = "'6.474"
$ws.Range("G6").Value = "'18"
$ws.Range("D7").Value = "'0.8135"
$ws.Range("G7").Value = "'18"
$ws.Range("D8").Value = "'0.8452"
$ws.Range("G8").Value = "'18"
$ws.Range("B9").Value = "BitrueCoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D9").Value = "'0.02845"
$ws.Range("E9").Value = "8BitrueCoinBTR"
$ws.Range("G9").Value = "'18"
$ws.Range("B10").Value = "BitMartToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D10").Value = "'0.09381"
$ws.Range("E10").Value = "9BitMartTokenBMX"
$ws.Range("G10").Value = "'18"
$ws.Range("B11").Value = "BitForexToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D11").Value = "'0.001510"
$ws.Range("E11").Value = "10BitForexTokenBF"
$ws.Range("G11").Value = "'18"
$ws.Range("B12").Value = "One"
$ws.Range("C12").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D12").Value = "'0.0005981"
$ws.Range("E12").Value = "11OneONE"
$ws.Range("G12").Value = "'18"
$ws.Range("B13").Value = "TigerCash"
$ws.Range("C13").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D13").Value = "'0.006123"
$ws.Range("E13").Value = "12TigerCashTCH"
$ws.Range("G13").Value = "'18"
$ws.Range("B14").Value = "LEO"
$ws.Range("C14").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D14").Value = "'3.600"
$ws.Range("E14").Value = "13LEOLEO"
$ws.Range("G14").Value = "'18"
$ws.Range("B15").Value = "GateToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D15").Value = "'3.010"
$ws.Range("E15").Value = "14GateTokenGT"
$ws.Range("G15").Value = "'18"
$ws.Range("B16").Value = "BTSEToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D16").Value = "'2.055"
$ws.Range("E16").Value = "15BTSETokenBTSE"
$ws.Range("G16").Value = "'18"
$ws.Range("B17").Value = "BitpandaEcosystemToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D17").Value = "'0.3207"
$ws.Range("E17").Value = "16BitpandaEcosystemTokenBEST"
$ws.Range("G17").Value = "'18"
$ws.Range("B18").Value = "WazirX"
$ws.Range("C18").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D18").Value = "'0.1339"
$ws.Range("E18").Value = "17WazirXWRX"
$ws.Range("G18").Value = "'18"
$ws.Range("D19").Value = "'0.06956"
$ws.Range("G19").Value = "'18"
$ws.Range("D20").Value = "'0.03195"
$ws.Range("G20").Value = "'18"
$ws.Range("D21").Value = "'0.1282"
$ws.Range("G21").Value = "'18"
$ws.Range("D22").Value = "'3.754"
$ws.Range("G22").Value = "'18"
$ws.Range("D23").Value = "'0.04651"
$ws.Range("G23").Value = "'18"
$ws.Range("G24").Value = "'18"
$ws.Range("D25").Value = "'0.001249"
$ws.Range("G25").Value = "'18"
$ws.Range("G26").Value = "'18"
$ws.Range("D27").Value = "'0.00009601"
$ws.Range("G27").Value = "'18"
$ws.Range("D28").Value = "'0.0001938"
$ws.Range("G28").Value = "'18"
$ws.Range("G29").Value = "'18"
$ws.Range("G30").Value = "'18"
$ws.Range("G31").Value = "'18"
$ws.Range("G32").Value = "'18"
$ws.Range("G33").Value = "'18"
$ws.Range("G34").Value = "'18"
$ws.Range("G35").Value = "'18"
$ws.Range("G36").Value = "'18"
$ws.Range("G37").Value = "'18"
$ws.Range("G38").Value = "'18"
$ws.Range("G39").Value = "'18"
$ws.Range("D40").Value = "'0.03668"
$ws.Range("G40").Value = "'18"
$ws.Range("D41").Value = "'0.006221"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("G41").Value = "'18"
$ws.Range("D42").Value = "'0.1056"
$ws.Range("G42").Value = "'18"
$ws.Range("D43").Value = "'0.002500"
$ws.Range("G43").Value = "'18"
$ws.Range("D44").Value = "'0.008746"
$ws.Range("G44").Value = "'18"
$ws.Range("D45").Value = "'0.00005295"
$ws.Range("G45").Value = "'18"
$ws.Range("G46").Value = "'18"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("G47").Value = "'18"
$ws.Range("D48").Value = "'0.002726"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("G48").Value = "'18"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("G49").Value = "'18"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("G50").Value = "'18"
$ws.Range("G51").Value = "'18"
